{"js": "// Change \"greatly\" to \"remarkably\" in the closing summary sentence:\n//   \"Git versioning is a powerful tool that, when used correctly, greatly\n//    enhances collaboration and project management in software development.\"\n// -> \"... remarkably enhances collaboration ...\"\n//\n// Using body.search keeps the run's existing character formatting\n// (Times New Roman / kern / ligatures) untouched since insertText(..., replace)\n// only swaps the text of the matched range.\nconst searchResults = context.document.body.search(\"greatly\", {\n  matchCase: true,\n  matchWholeWord: true\n});\nsearchResults.load(\"text\");\nawait context.sync();\n\nif (searchResults.items.length === 0) {\n  throw new Error('Could not find the word \"greatly\" to replace.');\n}\n\nsearchResults.items[0].insertText(\"remarkably\", Word.InsertLocation.replace);\nawait context.sync();\n", "ps1": "# Change \"greatly\" to \"remarkably\" in the closing summary sentence:\n#   \"Git versioning is a powerful tool that, when used correctly, greatly\n#    enhances collaboration and project management in software development.\"\n# -> \"... remarkably enhances collaboration ...\"\n#\n# Find & Replace on the document body preserves the run's existing character\n# formatting (Times New Roman / kern / ligatures) since only the matched\n# text is swapped out.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"greatly\"\n$find.MatchCase = $true\n$find.MatchWholeWord = $true\n$find.Replacement.Text = \"remarkably\"\n\n$wdReplaceOne = 1\n$wdFindContinue = 1\n\n$find.Execute(\n    [ref]$find.Text,        # FindText\n    [ref]$find.MatchCase,   # MatchCase\n    [ref]$find.MatchWholeWord, # MatchWholeWord\n    [ref]$false,            # MatchWildcards\n    [ref]$false,            # MatchSoundsLike\n    [ref]$false,            # MatchAllWordForms\n    [ref]$true,             # Forward\n    [ref]$wdFindContinue,   # Wrap\n    [ref]$false,            # Format\n    [ref]$find.Replacement.Text, # ReplaceWith\n    [ref]$wdReplaceOne      # Replace\n) | Out-Null\n"}
